$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.08
$ws.Cells.Item(3, 2).Value = 0.01668
$ws.Cells.Item(4, 2).Value = 0.04274
$ws.Cells.Item(5, 2).Value = 0.01605
$ws.Cells.Item(6, 2).Value = 0.03022
$ws.Cells.Item(7, 2).Value = 0.08111
$ws.Cells.Item(9, 2).Value = 0.00325
$ws.Cells.Item(10, 2).Value = 0.01008
$ws.Cells.Item(11, 2).Value = 0.01832
$ws.Cells.Item(12, 2).Value = 0.06494999999999999
$ws.Cells.Item(13, 2).Value = 0.0197
$ws.Cells.Item(14, 2).Value = 0.01026
$ws.Cells.Item(15, 2).Value = 0.03395
$ws.Cells.Item(16, 2).Value = 0.05152
$ws.Cells.Item(17, 2).Value = 0.03218
$ws.Cells.Item(18, 2).Value = 0.06619
$ws.Cells.Item(19, 2).Value = 0.10849
$ws.Cells.Item(20, 2).Value = 0.03113
$ws.Cells.Item(21, 2).Value = 0.04493
$ws.Cells.Item(22, 2).Value = 0.0531
$ws.Cells.Item(23, 2).Value = 0.06219
$ws.Cells.Item(24, 2).Value = 0.03046
$ws.Cells.Item(25, 2).Value = 0.0016
$ws.Cells.Item(26, 2).Value = 0.00851
$ws.Cells.Item(27, 2).Value = 0.00341
$ws.Cells.Item(28, 2).Value = 0.01434
$ws.Cells.Item(29, 2).Value = 0.00856
$ws.Cells.Item(30, 2).Value = 0.00342
$ws.Cells.Item(31, 2).Value = 0.00023
$ws.Cells.Item(32, 2).Value = 0.02025
$ws.Cells.Item(33, 2).Value = 0.00583
$ws.Cells.Item(34, 2).Value = 0.02633
